# abstraccion-links.xlsx — add "location" and "Move" HubSpot mapping
# tables below the existing content, and drop the now-redundant
# "HubSpot Contact" column (E) from the Move-properties table (rows 34-39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122
$xlVAlignBottom = -4107
$xlHAlignCenter = -4108

# ---------------------------------------------------------------------
# 1) Remove column E from the existing "Move property" table (rows 34-39)
#    Using Clear() removes the cell entirely (value + style), matching
#    the target which simply drops the <c> elements for E34:E39.
# ---------------------------------------------------------------------
$ws.Range("E34:E39").Clear() | Out-Null

# ---------------------------------------------------------------------
# 2) Seed the brand-new shared strings in the exact order they were
#    first typed by the original author so the shared string table
#    indices line up with the target workbook.
# ---------------------------------------------------------------------
$ws.Range("B44").Value = "propiedades locales y hubspot"
$ws.Range("B53").Value = "id_location"
$ws.Range("B52").Value = "numbre_areas"
$ws.Range("E52").Value = "number_of_areas"
$ws.Range("E54").Value = "country"
$ws.Range("B50").Value = "variable local"
$ws.Range("E50").Value = "variable Hubspot"
$ws.Range("E53").Value = "phone"
$ws.Range("B49").Value = "location"
$ws.Range("B58").Value = "Move"
$ws.Range("B60").Value = "id_move"

# ---------------------------------------------------------------------
# 3) Build five brand-new cell-formats as single "anchor" cells, in the
#    exact order they first appear so the generated cellXfs indices
#    line up with the target (11, 12, 13, 14, 15). The formatting is
#    then fanned out to every other cell that needs it via
#    Copy + PasteSpecial(Formats), which is the only reliable way this
#    runtime de-duplicates styles.
# ---------------------------------------------------------------------

# index 11: thin box border, no alignment override (table value cells)
$ws.Range("B51").Style = "Normal"
$ws.Range("B51").Borders.LineStyle = 1
$ws.Range("B51").Borders.Weight = 2

# index 12: bold font + thin box border + wrap text only (table headers)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B50").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B50").VerticalAlignment = $xlVAlignBottom

# index 13: big bold banner heading, centered, no border
$ws.Range("B44").Font.Bold = $true
$ws.Range("B44").Font.Size = 36
$ws.Range("B44").HorizontalAlignment = $xlHAlignCenter

# index 14: centered, no border (section sub-heading without rule)
$ws.Range("B58").Style = "Normal"
$ws.Range("B58").HorizontalAlignment = $xlHAlignCenter

# index 15: centered, thin bottom border only (section sub-heading with rule)
$ws.Range("B49").Style = "Normal"
$ws.Range("B49").HorizontalAlignment = $xlHAlignCenter
$ws.Range("B49").Borders.Item(9).LineStyle = 1
$ws.Range("B49").Borders.Item(9).Weight = 2

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Fan formatting out from the anchor cells to the rest of the table.
# ---------------------------------------------------------------------

# "propiedades locales y hubspot" banner -> whole B44:I47 block
$ws.Range("B44").Copy() | Out-Null
$ws.Range("B44:I47").PasteSpecial($xlPasteFormats) | Out-Null

# "location" sub-heading -> C49:E49
$ws.Range("B49").Copy() | Out-Null
$ws.Range("C49:E49").PasteSpecial($xlPasteFormats) | Out-Null

# "Move" sub-heading -> C58:E58
$ws.Range("B58").Copy() | Out-Null
$ws.Range("C58:E58").PasteSpecial($xlPasteFormats) | Out-Null

# table-header style (index 12) -> E50, B59, E59
$ws.Range("B50").Copy() | Out-Null
$ws.Range("E50").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B59").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E59").PasteSpecial($xlPasteFormats) | Out-Null

# existing "Property type" header style (index 1) -> C50:D50, C59:D59
$ws.Range("B5").Copy() | Out-Null
$ws.Range("C50:D50").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C59:D59").PasteSpecial($xlPasteFormats) | Out-Null

# table-value style (index 11) -> B/E columns of the data rows
$ws.Range("B51").Copy() | Out-Null
foreach ($row in 52..55) {
    $ws.Range("B$row").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("E$row").PasteSpecial($xlPasteFormats) | Out-Null
}
$ws.Range("E51").PasteSpecial($xlPasteFormats) | Out-Null
foreach ($row in 60..63) {
    $ws.Range("B$row").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("E$row").PasteSpecial($xlPasteFormats) | Out-Null
}

# existing "HubSpot Contact" value style (index 2) -> C/D columns of the data rows
$ws.Range("B6").Copy() | Out-Null
foreach ($row in 51..55) {
    $ws.Range("C$row`:D$row").PasteSpecial($xlPasteFormats) | Out-Null
}
foreach ($row in 60..63) {
    $ws.Range("C$row`:D$row").PasteSpecial($xlPasteFormats) | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5) Fill in the remaining cell values (existing shared strings, so
#    their write order has no effect on the shared string table).
# ---------------------------------------------------------------------

# "location" table (rows 49-55)
$ws.Range("C50").Value = "Location properties"
$ws.Range("D50").Value = "Property type"

$ws.Range("B51").Value = "name"
$ws.Range("C51").Value = "name"
$ws.Range("D51").Value = "string"
$ws.Range("E51").Value = "name"

$ws.Range("C52").Value = "Number of areas"
$ws.Range("D52").Value = "number"

$ws.Range("C53").Value = "id"
$ws.Range("D53").Value = "number"

$ws.Range("B54").Value = "region"
$ws.Range("C54").Value = "region"
$ws.Range("D54").Value = "string"

$ws.Range("B55").Value = "generation"
$ws.Range("C55").Value = "generation"
$ws.Range("D55").Value = "string"
$ws.Range("E55").Value = "generation"

# "Move" table (rows 58-63)
$ws.Range("B59").Value = "variable local"
$ws.Range("C59").Value = "Move property"
$ws.Range("D59").Value = "Object property"
$ws.Range("E59").Value = "variable Hubspot"

$ws.Range("C60").Value = "id"
$ws.Range("D60").Value = "Move Id"
$ws.Range("E60").Value = "id"

$ws.Range("B61").Value = "name"
$ws.Range("C61").Value = "name"
$ws.Range("D61").Value = "Name"
$ws.Range("E61").Value = "name"

$ws.Range("B62").Value = "pp"
$ws.Range("C62").Value = "PP"
$ws.Range("D62").Value = "PP"
$ws.Range("E62").Value = "pp"

$ws.Range("B63").Value = "power"
$ws.Range("C63").Value = "Power"
$ws.Range("D63").Value = "Power"
$ws.Range("E63").Value = "power"

# ---------------------------------------------------------------------
# 6) Merge the banner / section-heading cells.
# ---------------------------------------------------------------------
$ws.Range("B44:I47").Merge() | Out-Null
$ws.Range("B49:E49").Merge() | Out-Null
$ws.Range("B58:E58").Merge() | Out-Null

# ---------------------------------------------------------------------
# 7) Update the visible selection to reflect where the author ended up.
# ---------------------------------------------------------------------
$ws.Range("G78").Select() | Out-Null
